$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1692744.504728717
$ws.Range("D2").Value = 1691864.841200692
$ws.Range("E2").Value = 1692749.092345976
$ws.Range("F2").Value = 1691618.081491214
$ws.Range("G2").Value = 1691967.686920107
$ws.Range("H2").Value = 1693300.129814962
$ws.Range("I2").Value = 1691570.564286089
$ws.Range("J2").Value = 1692426.165228999
$ws.Range("K2").Value = 1692237.904945076
$ws.Range("L2").Value = 1692090.165077803
